$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused "Unnamed: 0.X" columns (old B:V) by deleting them and
# shifting the remaining columns (old W:Z -> new B:E) left.
$ws.Range("B1:V2").Delete(-4159) | Out-Null

# Update header row (row 1)
$ws.Range("B1").Value = "our_identified"
$ws.Range("C1").Value = "our_Overlap_merlin"
$ws.Range("D1").Value = "ourbest_param"
$ws.Range("E1").Value = "ourtime_taken"
$ws.Range("F1").Value = "best_paramcluster"
$ws.Range("G1").Value = "best_paramtraining"
$ws.Range("H1").Value = "best_paramwindow"
$ws.Range("I1").Value = "best_paramthreshold"

# New header cells (F1:I1) need the same bold/bordered header style as B1:E1
$ws.Range("B1").Copy() | Out-Null
$ws.Range("F1:I1").PasteSpecial(-4122) | Out-Null

# Update data row (row 2)
$ws.Range("B2").Value = "[]"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "{'cluster': 26, 'training': 386, 'window': 362, 'threshold': 1.5}"
$ws.Range("E2").Value = 20.03201633800199
$ws.Range("F2").Value = "RAS"
$ws.Range("G2").Value = "RAS"
$ws.Range("H2").Value = "RAS"
$ws.Range("I2").Value = "RAS"
